# Update "想去人数" (want-to-go count) figures in column F across the four
# sheets of the workbook, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1264
$ws.Range("F7").Value = 962
$ws.Range("F8").Value = 930
$ws.Range("F10").Value = 92
$ws.Range("F11").Value = 96
$ws.Range("F13").Value = 902
$ws.Range("F15").Value = 3758
$ws.Range("F16").Value = 1124
$ws.Range("F17").Value = 105
$ws.Range("F18").Value = 2537
$ws.Range("F20").Value = 1061
$ws.Range("F21").Value = 3498
$ws.Range("F22").Value = 730
$ws.Range("F24").Value = 35
$ws.Range("F25").Value = 2166
$ws.Range("F27").Value = 818
$ws.Range("F29").Value = 176
$ws.Range("F30").Value = 181
$ws.Range("F32").Value = 1295
$ws.Range("F33").Value = 1915
$ws.Range("F34").Value = 474
$ws.Range("F39").Value = 226

# Sheet 2: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 24

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 306

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1264
$ws.Range("F5").Value = 962
$ws.Range("F6").Value = 930
$ws.Range("F9").Value = 92
$ws.Range("F10").Value = 24
$ws.Range("F11").Value = 24
$ws.Range("F12").Value = 96
$ws.Range("F14").Value = 902
$ws.Range("F16").Value = 3758
$ws.Range("F17").Value = 1124
$ws.Range("F18").Value = 105
$ws.Range("F20").Value = 2537
$ws.Range("F22").Value = 1061
$ws.Range("F23").Value = 3498
$ws.Range("F24").Value = 730
$ws.Range("F27").Value = 35
$ws.Range("F28").Value = 2167
$ws.Range("F34").Value = 818
$ws.Range("F36").Value = 176
$ws.Range("F37").Value = 181
$ws.Range("F40").Value = 1295
$ws.Range("F41").Value = 1915
$ws.Range("F44").Value = 474
$ws.Range("F48").Value = 226
